$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix B1613 value (46030 -> 46029)
$ws.Cells.Item(1613, 2).Value2 = 46029

# 2. Defined name _FilterDatabase spans the (pre-existing) table.
# Must run before the new rows are appended below, otherwise the filter/defined-
# name ranges would be recalculated against the newly extended used range.
$n = $wb.Names.Item(1)
$n.RefersTo = '=in!$A$1:$L$1691'

# 3. Turn on the AutoFilter over the table range (also before extending the data,
# for the same reason as above).
$ws.Range("A1:L1691").AutoFilter() | Out-Null

# 4. New data rows 1679-1704.
# First stamp the formatting of the existing fully-populated row (1678) onto every
# target row so that styles (date format on column B, font on column A, etc.) match
# the rest of the table exactly.
for ($r = 1679; $r -le 1704; $r++) {
    $ws.Range("A1678:H1678").Copy($ws.Range("A$r`:H$r"))
}

$rows = @(
    @{R=1679; A="Alura";          B=46037; C="Banco Safra";       D=0;                   E=0;          F=0;              G=986064.78;          H=986064.78},
    @{R=1680; A="Alura";          B=46037; C="ITAU";               D=12024.65;            E=703317.9;   F=-691293.25;     G=361478.95;          H=-329814.3},
    @{R=1681; A="Alura";          B=46037; C="SANTANDER";          D=123214.12;           E=0;          F=123214.12;      G=0;                  H=123214.12},
    @{R=1682; A="Alura";          B=46037; C="XP investimentos "; D=0;                   E=0;          F=0;              G=26407046.870000001; H=26407046.870000001},
    @{R=1683; A="Alura";          B=46037; C="BANCO DO BRASIL";    D=0;                   E=0;          F=0;              G=39081.78;           H=39081.78},
    @{R=1684; A="FIAP";           B=46037; C="Banco Safra";       D=0;                   E=0;          F=0;              G=2749964.13;         H=2749964.13},
    @{R=1685; A="FIAP";           B=46037; C="ITAU";               D=36975.699999999997;  E=1651982.9;  F=-1615007.2;     G=390053.52;          H=-1224953.68},
    @{R=1686; A="FIAP";           B=46037; C="SANTANDER";          D=0;                   E=8052.91;    F=-8052.91;       G=59789.09;           H=51736.179999999993},
    @{R=1687; A="FIAP";           B=46037; C="XP investimentos "; D=0;                   E=0;          F=0;              G=51575007.030000001; H=51575007.030000001},
    @{R=1688; A="INSTITUTO FIAP"; B=46037; C="ITAU";               D=0;                   E=0;          F=0;              G=367835.82;          H=367835.82},
    @{R=1689; A="INSTITUTO FIAP"; B=46037; C="XP investimentos "; D=0;                   E=0;          F=0;              G=1156411.93;         H=1156411.93},
    @{R=1690; A="PM3";            B=46037; C="ITAU";               D=5551.24;             E=106410.11;  F=-100858.87;     G=87084.73;           H=-13774.14},
    @{R=1691; A="PM3";            B=46037; C="XP investimentos "; D=0;                   E=0;          F=0;              G=3591939.12;         H=3591939.12},
    @{R=1692; A="Alura";          B=46038; C="Banco Safra";       D=0;                   E=0;          F=0;              G=986064.78;          H=986064.78},
    @{R=1693; A="Alura";          B=46038; C="ITAU";               D=14664.12;            E=352969.64;  F=-338305.52;     G=108874.4;           H=-229431.12000000002},
    @{R=1694; A="Alura";          B=46038; C="SANTANDER";          D=164.26;              E=0;          F=164.26;         G=0;                  H=164.26},
    @{R=1695; A="Alura";          B=46038; C="XP investimentos "; D=0;                   E=0;          F=0;              G=26407046.870000001; H=26407046.870000001},
    @{R=1696; A="Alura";          B=46038; C="BANCO DO BRASIL";    D=0;                   E=0;          F=0;              G=31822.16;           H=31822.16},
    @{R=1697; A="FIAP";           B=46038; C="Banco Safra";       D=0;                   E=0;          F=0;              G=2749964.13;         H=2749964.13},
    @{R=1698; A="FIAP";           B=46038; C="ITAU";               D=92221.2;             E=0;          F=92221.2;        G=236832.48;          H=329053.68},
    @{R=1699; A="FIAP";           B=46038; C="SANTANDER";          D=0;                   E=12999.51;   F=-12999.51;      G=51736.18;           H=38736.67},
    @{R=1700; A="FIAP";           B=46038; C="XP investimentos "; D=0;                   E=0;          F=0;              G=51304899.509999998; H=51304899.509999998},
    @{R=1701; A="INSTITUTO FIAP"; B=46038; C="ITAU";               D=0;                   E=0;          F=0;              G=367835.82;          H=367835.82},
    @{R=1702; A="INSTITUTO FIAP"; B=46038; C="XP investimentos "; D=0;                   E=0;          F=0;              G=1158025.28;         H=1158025.28},
    @{R=1703; A="PM3";            B=46038; C="ITAU";               D=3854.61;             E=0;          F=3854.61;        G=0;                  H=3854.61},
    @{R=1704; A="PM3";            B=46038; C="XP investimentos "; D=0;                   E=0;          F=0;              G=3594093.95;         H=3594093.95}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $ws.Cells.Item($r, 5).Value2 = $row.E
    $ws.Cells.Item($r, 6).Value2 = $row.F
    $ws.Cells.Item($r, 7).Value2 = $row.G
    $ws.Cells.Item($r, 8).Value2 = $row.H
}

# 5. Update the view: scroll so the frozen pane shows row 1668, with C1668 selected.
$ws.Range("C1668").Select()
$excel.ActiveWindow.ScrollRow = 1668
$excel.ActiveWindow.ScrollColumn = 1
